$wb = $excel.ActiveWorkbook

# The "Swiss" sheet is the template for the new "Portugal" sheet: copy it
# to the end of the tab strip so it becomes the new last sheet.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Update the market-specific values on the new sheet. B4 is written before
# B2 so the new shared-string entries land in the same order as the target
# workbook ("NGC-3479/T2404" before "Portugal Market").
$portugal.Range("B4").Value = "NGC-3479/T2404"
$portugal.Range("B2").Value = "Portugal Market"

# Match the new sheet's column widths to the target layout.
$portugal.Columns.Item(1).ColumnWidth = 23.166666666666668
$portugal.Columns.Item(2).ColumnWidth = 15.166666666666666
$portugal.Columns.Item(3).ColumnWidth = 12.833333333333332
$portugal.Columns.Item(4).ColumnWidth = 14.666666666666666

# With the narrower layout, rows 3-5 wrap to a second line, doubling their height.
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Selection on the new sheet is just B2 (not the B2:B4 block inherited from Swiss).
$portugal.Range("B2").Select()

# The new Portugal sheet becomes the active tab.
$portugal.Activate()
